# Update "想去人数" (want-to-go count) figures in column F across sheets,
# matching the regenerated data output at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 880
$ws1.Range("F6").Value = 499
$ws1.Range("F8").Value = 223
$ws1.Range("F10").Value = 69
$ws1.Range("F11").Value = 547
$ws1.Range("F12").Value = 0

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0
$ws2.Range("F4").Value = 0

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 0
$ws4.Range("F4").Value = 337
$ws4.Range("F5").Value = 0
$ws4.Range("F7").Value = 4425
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 223
$ws4.Range("F12").Value = 498
$ws4.Range("F13").Value = 69
$ws4.Range("F14").Value = 0
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 547
